# Deploy versión estable desde local
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old header row (DIV 1 / DIV 2 / DIV 3 - A / DIV 3 - B); this
# shifts everything up by one row, matching columns A and B exactly.
$ws.Rows(1).Delete()

# Column C (division 3 - A) and column D (division 3 - B) are re-entered
# with an updated order, and one roster name is corrected.
$colC = @("Miguel Jumbo", "Jorge Diaz", "David Tamarit", "Eudes Vázquez", "Iñaki Esnal", "María Toral", "Álvaro Monleon")
$colD = @("Enrique Lázaro", "Antonio Calatayud", "Attila Bajner", "Leo Garrido", "Alejandro Navarro", "Miguel Mau", "Jorge Alberola")

for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $colD[$i]
}

$ws.Range("A1:XFD1").Select()
